$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.286832544864788
$ws.Range("C2").Value = 1.655778082260271
$ws.Range("D2").Value = 0.7527432677738641
$ws.Range("E2").Value = 0.4942365360607697
$ws.Range("G2").Value = 6.189590430959694

$ws.Range("B3").Value = 1.455362044514542
$ws.Range("C3").Value = 1.655778082260271
$ws.Range("D3").Value = 3.537761648806719
$ws.Range("E3").Value = 10.19245300693656
$ws.Range("G3").Value = 16.84135478251809
